$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Förändrad" (column C) date for rows 2-16 from 2023-10-22 to 2023-10-25
$ws.Range("C2:C16").Value = 45224
